# Applies the "Corrections 0903 (suite réunion 0803)" edit to
# StructureDefinition-ror-practitioner.xlsx
#
# 1. Metadata sheet: bump the "Date" property value.
# 2. Elements sheet:
#    - New French "Short" descriptions for several rows
#      (identifier, telecom, telecom.value).
#    - Re-order the three Practitioner.telecom.extension slice rows
#      (usage / confidentiality-level / communication-channel) and give
#      each of them a new French "Short" description.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet - "Date" row (row 8, column B)
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2023-03-09T09:15:18+00:00"

# ---------------------------------------------------------------------
# 2. Elements sheet
# ---------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Column map used below (1-based):
#   A=1 ID, C=3 Slice Name, F=6 Min, K=11 Type(s), L=12 Short, M=13 Definition

# Practitioner.identifier (row 20) - new Short text
$els.Cells.Item(20, 12).Value = "idNat_PS (Professionnel) : Identification nationale du professionnel définie par le CI-SIS"

# Practitioner.telecom (row 23) - new Short text
$els.Cells.Item(23, 12).Value = "boiteLettreMSS (Professionnel) : Boîte(s) aux lettres du service de messagerie sécurisée de santé (MSS) rattachée(s) au professionnel"

# Practitioner.telecom.extension slices (rows 26-28) are reordered so the
# "communication-channel" slice now comes first, followed by "usage" and
# then "confidentiality-level". Each also receives a new French Short text.

# Row 26 -> ror-telecom-communication-channel
$els.Cells.Item(26, 1).Value  = "Practitioner.telecom.extension:ror-telecom-communication-channel"
$els.Cells.Item(26, 3).Value  = "ror-telecom-communication-channel"
$els.Cells.Item(26, 6).Value  = "'1"
$els.Cells.Item(26, 11).Value = "Extension {https://interop.esante.gouv.fr/ig/fhir/ror30/StructureDefinition/ror-telecom-communication-channel}`n"
$els.Cells.Item(26, 12).Value = "canal (Telecommunication) : Code spécifiant le canal ou la manière dont s'établit la communication"
$els.Cells.Item(26, 13).Value = "Extension créée dans le cadre du ROR spécifiant le canal ou la manière dont s'établit la communication "

# Row 27 -> ror-telecom-usage
$els.Cells.Item(27, 1).Value  = "Practitioner.telecom.extension:ror-telecom-usage"
$els.Cells.Item(27, 3).Value  = "ror-telecom-usage"
$els.Cells.Item(27, 6).Value  = "'0"
$els.Cells.Item(27, 11).Value = "Extension {https://interop.esante.gouv.fr/ig/fhir/ror30/StructureDefinition/ror-telecom-usage}`n"
$els.Cells.Item(27, 12).Value = "utilisation (Telecommunication) : Utilisation du canal de communication "
$els.Cells.Item(27, 13).Value = "Extension créée dans le cadre du ROR qui précise l'utilisation du canal de communication "

# Row 28 -> ror-telecom-confidentiality-level
$els.Cells.Item(28, 1).Value  = "Practitioner.telecom.extension:ror-telecom-confidentiality-level"
$els.Cells.Item(28, 3).Value  = "ror-telecom-confidentiality-level"
$els.Cells.Item(28, 6).Value  = "1"
$els.Cells.Item(28, 11).Value = "Extension {https://interop.esante.gouv.fr/ig/fhir/ror30/StructureDefinition/ror-telecom-confidentiality-level}`n"
$els.Cells.Item(28, 12).Value = "niveauConfidentialite (Telecommunication) : niveau de restriction de l'accès aux attributs de la classe Télécommunication"
$els.Cells.Item(28, 13).Value = "Extension créée dans le cadre du ROR qui permet de définir le niveau de restriction de l'accès aux attributs de la classe Télécommunication."

# Practitioner.telecom.value (row 30) - new Short text
$els.Cells.Item(30, 12).Value = "adresseTelecom (Telecommunication) : Valeur de l'adresse de télécommunication dans le format induit par le canal de communication"
